{"js": "// Replace the two-digit \u00f7 one-digit division problems in the table with\n// the new set of problems, preserving cell formatting (font/size/etc.)\n// by doing an in-place text replacement via Range.insertText(..., \"Replace\").\nconst replacements = [\n  [\"98\u00f76=\", \"27\u00f77=\"],\n  [\"33\u00f73=\", \"75\u00f75=\"],\n  [\"65\u00f73=\", \"79\u00f76=\"],\n  [\"50\u00f75=\", \"62\u00f74=\"],\n  [\"32\u00f75=\", \"23\u00f79=\"],\n  [\"71\u00f76=\", \"96\u00f74=\"],\n  [\"72\u00f73=\", \"94\u00f75=\"],\n  [\"19\u00f77=\", \"95\u00f72=\"],\n  [\"63\u00f78=\", \"22\u00f75=\"],\n  [\"30\u00f72=\", \"78\u00f75=\"],\n  [\"50\u00f78=\", \"92\u00f78=\"],\n  [\"30\u00f78=\", \"65\u00f74=\"],\n  [\"66\u00f72=\", \"84\u00f77=\"],\n  [\"28\u00f77=\", \"24\u00f77=\"],\n  [\"21\u00f75=\", \"24\u00f73=\"],\n  [\"79\u00f72=\", \"53\u00f74=\"],\n  [\"24\u00f78=\", \"50\u00f73=\"],\n  [\"61\u00f74=\", \"92\u00f75=\"],\n  [\"28\u00f73=\", \"83\u00f77=\"],\n  [\"20\u00f73=\", \"81\u00f78=\"],\n  [\"29\u00f72=\", \"21\u00f75=\"],\n  [\"87\u00f74=\", \"50\u00f78=\"],\n  [\"45\u00f75=\", \"49\u00f75=\"],\n  [\"83\u00f76=\", \"48\u00f79=\"],\n  [\"67\u00f72=\", \"31\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit \u00f7 one-digit division problems in the table with\n# the new set of problems, preserving cell formatting by using Find/Replace\n# on the document's content range (each old value is unique in the doc).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"98\u00f76=\", \"27\u00f77=\"),\n    @(\"33\u00f73=\", \"75\u00f75=\"),\n    @(\"65\u00f73=\", \"79\u00f76=\"),\n    @(\"50\u00f75=\", \"62\u00f74=\"),\n    @(\"32\u00f75=\", \"23\u00f79=\"),\n    @(\"71\u00f76=\", \"96\u00f74=\"),\n    @(\"72\u00f73=\", \"94\u00f75=\"),\n    @(\"19\u00f77=\", \"95\u00f72=\"),\n    @(\"63\u00f78=\", \"22\u00f75=\"),\n    @(\"30\u00f72=\", \"78\u00f75=\"),\n    @(\"50\u00f78=\", \"92\u00f78=\"),\n    @(\"30\u00f78=\", \"65\u00f74=\"),\n    @(\"66\u00f72=\", \"84\u00f77=\"),\n    @(\"28\u00f77=\", \"24\u00f77=\"),\n    @(\"21\u00f75=\", \"24\u00f73=\"),\n    @(\"79\u00f72=\", \"53\u00f74=\"),\n    @(\"24\u00f78=\", \"50\u00f73=\"),\n    @(\"61\u00f74=\", \"92\u00f75=\"),\n    @(\"28\u00f73=\", \"83\u00f77=\"),\n    @(\"20\u00f73=\", \"81\u00f78=\"),\n    @(\"29\u00f72=\", \"21\u00f75=\"),\n    @(\"87\u00f74=\", \"50\u00f78=\"),\n    @(\"45\u00f75=\", \"49\u00f75=\"),\n    @(\"83\u00f76=\", \"48\u00f79=\"),\n    @(\"67\u00f72=\", \"31\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
